# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" suffixed header columns to "_FV2404" / "_FV2410"
# 2) Turn the data range into an Excel Table (ListObject)
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 = "_old" -> "_FV2404", L1:U1 = "_new" -> "_FV2410"; K1 "diff" stays) ---
$headerRange = $ws.Range("A1:J1")
$null = $headerRange.Replace("_old", "_FV2404", 2)

$headerRange2 = $ws.Range("L1:U1")
$null = $headerRange2.Replace("_new", "_FV2410", 2)

# --- 2. Convert the used range into a native Excel table ---
$tableRange = $ws.Range("A1:U78")
$tbl = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$tbl.Name = "Table1"

# --- 3. Freeze panes at row 2 (so row 1 header stays visible) ---
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
